$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Tabelle1")
$ws.Activate()

# The author resized their Excel window before saving - reflect that on the
# active window (best effort; harmless if the host does not expose it).
$excel.ActiveWindow.Width = 13410
$excel.ActiveWindow.Height = 7515

# Row 24 (the "Calc" row) used to store the text "x" in the Numeric /
# Pseudonumeric / Pseudo-w-Unicode / String / MultiValue columns. It should
# hold the numeric value 1 instead, matching every other data row.
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("K24").Value = 1

# Scroll the frozen bottom-right pane so E15 is the new top-left visible
# cell, then leave K23 selected (both differ from the previous B15 / A24).
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K23").Select() | Out-Null
